$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

$values = @(
    "wm3pc5kieot",
    "n3bk8u16",
    "HHKD",
    "Madam Thu Bakery, 21C, Võ Văn Tần, Ninh Kiều, Ninh Kiều District, Cần Thơ, 94111, Vietnam",
    "https://www.google.com/maps/search/?api=1&query=10.032100,105.786400",
    "2025-08-21T07:29:45.242Z",
    "",
    "",
    "",
    "",
    "",
    "Nguyễn Văn B",
    "d7ee3a393285b163",
    "f8fa057026afe9e0606ed5c0c911b036239a6d9486af784c90ea1edb7325a99a"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}
